# Auto-applies the numeric corrections recorded for the Diabolos_Profits leve-profit
# workbook. Each worksheet corresponds to a crafting class (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# For every affected leve row we rewrite the recalculated currentAveragePrice /
# LevePrice / LeveProfit figures (columns H, I, J, K, L, M, N).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 233143.98
$ws.Range("J17").Value = 233143.98
$ws.Range("L17").Value = 699431.9400000001
$ws.Range("N17").Value = -699767.9400000001
# Row 125
$ws.Range("H125").Value = 1553.5
$ws.Range("I125").Value = 1221
$ws.Range("J125").Value = 1886
$ws.Range("K125").Value = 10989
$ws.Range("L125").Value = 16974
$ws.Range("M125").Value = -8529
$ws.Range("N125").Value = -21894
# Row 132
$ws.Range("H132").Value = 4294.724
$ws.Range("I132").Value = 4353.6294
$ws.Range("K132").Value = 13060.8882
$ws.Range("M132").Value = -10530.8882
# Row 138
$ws.Range("H138").Value = 2762.658
$ws.Range("I138").Value = 1223.4348
$ws.Range("J138").Value = 5122.8
$ws.Range("K138").Value = 3670.3044
$ws.Range("L138").Value = 15368.4
$ws.Range("M138").Value = 1469.6956
$ws.Range("N138").Value = -25648.4

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1153.9
$ws.Range("I2").Value = 1328.1428
$ws.Range("J2").Value = 747.3333
$ws.Range("K2").Value = 1328.1428
$ws.Range("L2").Value = 747.3333
$ws.Range("M2").Value = -1215.1428
$ws.Range("N2").Value = -973.3333
# Row 32
$ws.Range("H32").Value = 5017.405
$ws.Range("I32").Value = 3006.8823
$ws.Range("J32").Value = 13562.125
$ws.Range("K32").Value = 3006.8823
$ws.Range("L32").Value = 13562.125
$ws.Range("M32").Value = -2719.8823
$ws.Range("N32").Value = -14136.125
# Row 61
$ws.Range("H61").Value = 47621276
$ws.Range("I61").Value = 83334460
$ws.Range("J61").Value = 3697
$ws.Range("K61").Value = 83334460
$ws.Range("L61").Value = 3697
$ws.Range("M61").Value = -83334248
$ws.Range("N61").Value = -4121
# Row 63
$ws.Range("H63").Value = 70593700
$ws.Range("I63").Value = 100006136
$ws.Range("K63").Value = 100006136
$ws.Range("M63").Value = -100005450
# Row 66
$ws.Range("H66").Value = 70593700
$ws.Range("I66").Value = 100006136
$ws.Range("K66").Value = 500030680
$ws.Range("M66").Value = -500027248
# Row 74
$ws.Range("H74").Value = 27779868
$ws.Range("I74").Value = 33335046
$ws.Range("K74").Value = 33335046
$ws.Range("M74").Value = -33334172
# Row 77
$ws.Range("H77").Value = 27779868
$ws.Range("I77").Value = 33335046
$ws.Range("K77").Value = 166675230
$ws.Range("M77").Value = -166670862
# Row 108
$ws.Range("H108").Value = 41998.5
$ws.Range("J108").Value = 41998.5
$ws.Range("L108").Value = 41998.5
$ws.Range("N108").Value = -49678.5
# Row 116
$ws.Range("H116").Value = 1153.9
$ws.Range("I116").Value = 1328.1428
$ws.Range("J116").Value = 747.3333
$ws.Range("K116").Value = 1328.1428
$ws.Range("L116").Value = 747.3333
$ws.Range("M116").Value = 965.8571999999999
$ws.Range("N116").Value = -5335.3333
# Row 122
$ws.Range("H122").Value = 27779780
$ws.Range("I122").Value = 33335088
$ws.Range("K122").Value = 100005264
$ws.Range("M122").Value = -100002814
# Row 132
$ws.Range("H132").Value = 4200
$ws.Range("I132").Value = 3900
$ws.Range("K132").Value = 11700
$ws.Range("M132").Value = -9170
# Row 136
$ws.Range("H136").Value = 47621276
$ws.Range("I136").Value = 83334460
$ws.Range("J136").Value = 3697
$ws.Range("K136").Value = 250003380
$ws.Range("L136").Value = 11091
$ws.Range("M136").Value = -250000830
$ws.Range("N136").Value = -16191

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1153.9
$ws.Range("I3").Value = 1328.1428
$ws.Range("J3").Value = 747.3333
$ws.Range("K3").Value = 1328.1428
$ws.Range("L3").Value = 747.3333
$ws.Range("M3").Value = -1214.1428
$ws.Range("N3").Value = -975.3333
# Row 99
$ws.Range("H99").Value = 1655.6364
$ws.Range("I99").Value = 1141.4
$ws.Range("K99").Value = 1141.4
$ws.Range("M99").Value = 356.5999999999999
# Row 134
$ws.Range("H134").Value = 1933
$ws.Range("I134").Value = 1933
$ws.Range("K134").Value = 5799
$ws.Range("M134").Value = -3264

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3321.614
$ws.Range("I31").Value = 1791.75
$ws.Range("J31").Value = 3571.3877
$ws.Range("K31").Value = 1791.75
$ws.Range("L31").Value = 3571.3877
$ws.Range("M31").Value = -1496.75
$ws.Range("N31").Value = -4161.3877
# Row 34
$ws.Range("H34").Value = 3321.614
$ws.Range("I34").Value = 1791.75
$ws.Range("J34").Value = 3571.3877
$ws.Range("K34").Value = 1791.75
$ws.Range("L34").Value = 3571.3877
$ws.Range("M34").Value = -1589.75
$ws.Range("N34").Value = -3975.3877
# Row 58
$ws.Range("H58").Value = 1361.64
$ws.Range("I58").Value = 1162.75
$ws.Range("J58").Value = 2157.2
$ws.Range("K58").Value = 1162.75
$ws.Range("L58").Value = 2157.2
$ws.Range("M58").Value = -959.75
$ws.Range("N58").Value = -2563.2
# Row 102
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
# Row 134
$ws.Range("H134").Value = 1505.8667
$ws.Range("I134").Value = 1428.8462
$ws.Range("J134").Value = 2006.5
$ws.Range("K134").Value = 4286.5386
$ws.Range("L134").Value = 6019.5
$ws.Range("M134").Value = -1751.5386
$ws.Range("N134").Value = -11089.5
# Row 136
$ws.Range("H136").Value = 1361.64
$ws.Range("I136").Value = 1162.75
$ws.Range("J136").Value = 2157.2
$ws.Range("K136").Value = 3488.25
$ws.Range("L136").Value = 6471.599999999999
$ws.Range("M136").Value = -938.25
$ws.Range("N136").Value = -11571.6

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 947.1053000000001
$ws.Range("I5").Value = 732.9167
$ws.Range("K5").Value = 2198.7501
$ws.Range("M5").Value = -2086.7501
# Row 135
$ws.Range("H135").Value = 947.1053000000001
$ws.Range("I135").Value = 732.9167
$ws.Range("K135").Value = 6596.2503
$ws.Range("M135").Value = -4061.2503

$ws = $wb.Worksheets.Item("GSM")
# Row 58
$ws.Range("H58").Value = 23989.4
$ws.Range("I58").Value = 22499.5
$ws.Range("J58").Value = 24982.666
$ws.Range("K58").Value = 22499.5
$ws.Range("L58").Value = 24982.666
$ws.Range("M58").Value = -22222.5
$ws.Range("N58").Value = -25536.666
# Row 70
$ws.Range("H70").Value = 6344.2
$ws.Range("I70").Value = 7243.6875
$ws.Range("J70").Value = 4745.1113
$ws.Range("K70").Value = 7243.6875
$ws.Range("L70").Value = 4745.1113
$ws.Range("M70").Value = -6973.6875
$ws.Range("N70").Value = -5285.1113
# Row 73
$ws.Range("H73").Value = 6344.2
$ws.Range("I73").Value = 7243.6875
$ws.Range("J73").Value = 4745.1113
$ws.Range("K73").Value = 7243.6875
$ws.Range("L73").Value = 4745.1113
$ws.Range("M73").Value = -6307.6875
$ws.Range("N73").Value = -6617.1113
# Row 122
$ws.Range("H122").Value = 2889.2856
$ws.Range("I122").Value = 2413.4285
$ws.Range("J122").Value = 3365.1428
$ws.Range("K122").Value = 7240.2855
$ws.Range("L122").Value = 10095.4284
$ws.Range("M122").Value = -4790.2855
$ws.Range("N122").Value = -14995.4284
# Row 132
$ws.Range("H132").Value = 6036.263
$ws.Range("I132").Value = 5593.273
$ws.Range("J132").Value = 6645.375
$ws.Range("K132").Value = 16779.819
$ws.Range("L132").Value = 19936.125
$ws.Range("M132").Value = -14249.819
$ws.Range("N132").Value = -24996.125

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()
# Row 40
$ws.Range("H40").Value = 2785.7144
$ws.Range("J40").Value = 2000
$ws.Range("L40").Value = 2000
$ws.Range("N40").Value = -2272
# Row 100
$ws.Range("H100").Value = 3489.4546
$ws.Range("I100").Value = 3438.4
$ws.Range("J100").Value = 4000
$ws.Range("K100").Value = 3438.4
$ws.Range("L100").Value = 4000
$ws.Range("M100").Value = -2897.4
$ws.Range("N100").Value = -5082
# Row 110
$ws.Range("H110").Value = 14666.333
$ws.Range("J110").Value = 14666.333
$ws.Range("L110").Value = 14666.333
$ws.Range("N110").Value = -22846.333

$ws = $wb.Worksheets.Item("WVR")
# Row 105
$ws.Range("H105").Value = 30615
$ws.Range("J105").Value = 30615
$ws.Range("L105").Value = 30615
$ws.Range("N105").Value = -37603
